$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "60.944.00"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.650.20"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "577.51"
$ws.Range("E5").Value = "  +0.15%  "

Set-TextValue "D6" "144.49"
$ws.Range("E6").Value = "  +1.03%  "

Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("E10").Value = "  +0.71%  "

Set-TextValue "D11" "0.381"
$ws.Range("E11").Value = "  +3.47%  "

Set-TextValue "D12" "0.155"
$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("D13").Value = "3.125.28"
$ws.Range("E13").Value = "  +1.07%  "

Set-TextValue "D14" "26.12"
$ws.Range("E14").Value = "  +12.01%  "

$ws.Range("D15").Value = "60.983.31"
$ws.Range("E15").Value = "  +0.59%  "

Set-TextValue "D16" "0.0000143"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "2.667.19"
$ws.Range("E17").Value = "  +1.45%  "

Set-TextValue "D18" "11.63"
$ws.Range("E18").Value = "  +3.07%  "

Set-TextValue "D19" "4.72"
$ws.Range("E19").Value = "  +1.52%  "

Set-TextValue "D20" "349.75"
$ws.Range("E20").Value = "  +0.27%  "

Set-TextValue "D21" "6.94"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  +1.98%  "

Set-TextValue "D24" "64.13"
$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("E26").Value = "  -0.17%  "

Set-TextValue "D27" "8.15"
$ws.Range("E27").Value = "  +4.96%  "

Set-TextValue "D28" "1.97"
$ws.Range("E28").Value = "  +7.48%  "

$ws.Range("D29").Value = "0.0₃0811"
$ws.Range("E29").Value = "  +2.20%  "

Set-TextValue "D30" "6.84"
$ws.Range("E30").Value = "  +7.57%  "

$ws.Range("E31").Value = "  +0.04%  "

Set-TextValue "D32" "166.88"
$ws.Range("E32").Value = "  +2.69%  "

Set-TextValue "D33" "19.91"
$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("E34").Value = "  +6.86%  "

$ws.Range("E35").Value = "  +7.26%  "

$ws.Range("E36").Value = "  +7.40%  "

Set-TextValue "D37" "1.67"
$ws.Range("E37").Value = "  +3.80%  "

Set-TextValue "D38" "337.69"
$ws.Range("E38").Value = "  +12.34%  "

Set-TextValue "D39" "4.04"
$ws.Range("E39").Value = "  +3.87%  "

Set-TextValue "D40" "0.903"
$ws.Range("E40").Value = "  +6.64%  "

Set-TextValue "D41" "38.58"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("E42").Value = "  +4.83%  "

Set-TextValue "D43" "20.44"
$ws.Range("E43").Value = "  +2.56%  "

Set-TextValue "D44" "134.02"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.0996"
$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0561"
$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.616"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0248"
$ws.Range("E48").Value = "  +2.47%  "

Set-TextValue "D49" "20.52"
$ws.Range("E49").Value = "  +2.54%  "

$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").Value = "2.097.41"
$ws.Range("E51").Value = "  +3.61%  "
